$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content (old table used A1:J8)
$ws.Cells.Clear()

# Rewrite the table with the new columns/rows
$ws.Cells.Item(1, 1).Value = 'name'
$ws.Cells.Item(1, 2).Value = 'email'
$ws.Cells.Item(1, 3).Value = 'gender'
$ws.Cells.Item(1, 4).Value = 'address'
$ws.Cells.Item(1, 5).Value = 'age'
$ws.Cells.Item(1, 6).Value = 'role'

$ws.Cells.Item(2, 1).Value = 'I''m super admin'
$ws.Cells.Item(2, 2).Value = 'admin@gmail.com'
$ws.Cells.Item(2, 3).Value = 'MALE'
$ws.Cells.Item(2, 4).Value = 'hn'
$ws.Cells.Item(2, 5).Value = 25
$ws.Cells.Item(2, 6).Value = 1

$ws.Cells.Item(3, 1).Value = 'nu '
$ws.Cells.Item(3, 2).Value = 'nhokanhanh@gmail.com'
$ws.Cells.Item(3, 3).Value = 'MALE'
$ws.Cells.Item(3, 4).Value = 'pham ghong thai'
$ws.Cells.Item(3, 5).Value = 23
$ws.Cells.Item(3, 6).Value = 1

$ws.Cells.Item(4, 1).Value = 'hr'
$ws.Cells.Item(4, 2).Value = 'hr@gmail.com'
$ws.Cells.Item(4, 3).Value = 'MALE'
$ws.Cells.Item(4, 4).Value = '42 Phạm Hồng Thái'
$ws.Cells.Item(4, 5).Value = 25
$ws.Cells.Item(4, 6).Value = 1

$ws.Cells.Item(5, 1).Value = 'loe vann nguyen'
$ws.Cells.Item(5, 2).Value = 'test@gmail.com'
$ws.Cells.Item(5, 3).Value = 'MALE'
$ws.Cells.Item(5, 4).Value = '41b le van tho'
$ws.Cells.Item(5, 5).Value = 23
$ws.Cells.Item(5, 6).Value = 1

$ws.Cells.Item(6, 1).Value = 'sos nu'
$ws.Cells.Item(6, 2).Value = 'sosnu1111@gmail.com'
$ws.Cells.Item(6, 3).Value = 'MALE'
$ws.Cells.Item(6, 4).Value = 'ha-noi'
$ws.Cells.Item(6, 5).Value = 25
$ws.Cells.Item(6, 6).Value = 1

$ws.Cells.Item(7, 1).Value = 'sos nu'
$ws.Cells.Item(7, 2).Value = 'sosnu@gmail.com'
$ws.Cells.Item(7, 3).Value = 'MALE'
$ws.Cells.Item(7, 4).Value = 'ha-noi'
$ws.Cells.Item(7, 5).Value = 25
$ws.Cells.Item(7, 6).Value = 1

$ws.Cells.Item(8, 1).Value = 'sos nu1'
$ws.Cells.Item(8, 2).Value = 'sosnu11111111@gmail.com'
$ws.Cells.Item(8, 3).Value = 'MALE'
$ws.Cells.Item(8, 4).Value = 'ha-noi'
$ws.Cells.Item(8, 5).Value = 25
$ws.Cells.Item(8, 6).Value = 1

$ws.Cells.Item(9, 1).Value = 'sos nu'
$ws.Cells.Item(9, 2).Value = 'so11sn112313123u@gmail.com'
$ws.Cells.Item(9, 3).Value = 'MALE'
$ws.Cells.Item(9, 4).Value = 'ha-noi'
$ws.Cells.Item(9, 5).Value = 25
$ws.Cells.Item(9, 6).Value = 1

$ws.Cells.Item(10, 1).Value = 'sos nu1'
$ws.Cells.Item(10, 2).Value = 'sosnu111111231111@gmail.com'
$ws.Cells.Item(10, 3).Value = 'MALE'
$ws.Cells.Item(10, 4).Value = 'ha-noi'
$ws.Cells.Item(10, 5).Value = 25
$ws.Cells.Item(10, 6).Value = 1

$ws.Cells.Item(11, 1).Value = 'sos nu'
$ws.Cells.Item(11, 2).Value = 'so11sn112312313123u@gmail.com'
$ws.Cells.Item(11, 3).Value = 'MALE'
$ws.Cells.Item(11, 4).Value = 'ha-noi'
$ws.Cells.Item(11, 5).Value = 25
$ws.Cells.Item(11, 6).Value = 1

$ws.Cells.Item(12, 1).Value = 'sos nu1'
$ws.Cells.Item(12, 2).Value = 'sosnu11111111231111@gmail.com'
$ws.Cells.Item(12, 3).Value = 'MALE'
$ws.Cells.Item(12, 4).Value = 'ha-noi'
$ws.Cells.Item(12, 5).Value = 25
$ws.Cells.Item(12, 6).Value = 1

# Rename the sheet and tab
$ws.Name = 'export-user (6)'

# Restore the active selection to F11 as in the saved workbook
$ws.Range('F11').Select()
